$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.25 = 12534.15 pesos`n✅ 12534.15 pesos = 3.23 = 968.89 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 307.8
$wsTasas.Range("O10").Value = 3858.01
$wsTasas.Range("N12").Value = 3881
$wsTasas.Range("O12").Value = 300.001
